$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 8512
$ws.Range("F5").Value = 569
$ws.Range("F6").Value = 7539
$ws.Range("F12").Value = 203
$ws.Range("F13").Value = 12581
$ws.Range("F15").Value = 110
$ws.Range("F16").Value = 2720
$ws.Range("F18").Value = 5908
$ws.Range("F20").Value = 3127
$ws.Range("F24").Value = 25
$ws.Range("F28").Value = 3454
$ws.Range("F29").Value = 71
$ws.Range("F30").Value = 2559
$ws.Range("F32").Value = 1798
$ws.Range("F33").Value = 88
$ws.Range("F34").Value = 166
$ws.Range("F35").Value = 6255
$ws.Range("F37").Value = 149
$ws.Range("F38").Value = 1271
$ws.Range("F39").Value = 59
$ws.Range("F40").Value = 969
$ws.Range("F42").Value = 205
$ws.Range("F44").Value = 1124
$ws.Range("F46").Value = 1134
$ws.Range("F47").Value = 1652
$ws.Range("F48").Value = 38
$ws.Range("F50").Value = 1148

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 31
$ws.Range("F7").Value = 13
$ws.Range("F17").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 372
$ws.Range("F3").Value = 534

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 372
$ws.Range("F6").Value = 8512
$ws.Range("F7").Value = 569
$ws.Range("F8").Value = 7539
$ws.Range("F14").Value = 203
$ws.Range("F16").Value = 12581
$ws.Range("F18").Value = 2720
$ws.Range("F19").Value = 5908
$ws.Range("F20").Value = 3127
$ws.Range("F27").Value = 3455
$ws.Range("F28").Value = 71
$ws.Range("F29").Value = 2559
$ws.Range("F31").Value = 1798
$ws.Range("F32").Value = 166
$ws.Range("F33").Value = 6255
$ws.Range("F36").Value = 149
$ws.Range("F38").Value = 1271
$ws.Range("F39").Value = 59
$ws.Range("F40").Value = 969
$ws.Range("F42").Value = 205
$ws.Range("F44").Value = 1124
$ws.Range("F46").Value = 1134
$ws.Range("F47").Value = 1652
$ws.Range("F48").Value = 38
$ws.Range("F50").Value = 1148

